$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header and data values to focus on "Job Postings" instead of "Sales"
$ws.Range("B1").Value = "Job Postings"
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 15

# Auto-fit column B to match the new, narrower content width (stored width 12)
$ws.Columns.Item(2).ColumnWidth = 11.17

# Move the active selection to H7
$ws.Range("H7").Select()
